$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Typography": row 13 (Typography_06 / verdana.ttf) changes its size
# from 20 to 12, and a new row 14 is inserted for a new typography entry
# "Typography_07" that reuses verdana.ttf at size 20 (what row 13 used to be).
# ---------------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

$typo.Range("D13").Value = 12

$typo.Range("B14").Value = "Typography_07"
$typo.Range("C14").Value = "verdana.ttf"
$typo.Range("D14").Value = 20
$typo.Range("E14").Value = 4
$typo.Range("F14").Value = "?"
$typo.Range("G14").Value = ""
$typo.Range("H14").Value = ""
$typo.Range("I14").Value = ""
$typo.Range("J14").Value = ""

# ---------------------------------------------------------------------------
# Sheet "Translation": 86 new rows (19-104) describing per-digit steering
# wheel status / fuel-usage calculation parameter & value text entries.
# Every row alternates between a "PARAMETER" placeholder (Typography_06,
# left aligned, LTR) and a "val" placeholder (Typography_02, left aligned,
# LTR). Text ids follow "SingleUseId18".."SingleUseId103".
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

$row = 19
$singleUseId = 18
for ($pair = 0; $pair -lt 43; $pair++) {
    $tr.Range("B" + $row).Value = "SingleUseId" + $singleUseId
    $tr.Range("C" + $row).Value = "Typography_06"
    $tr.Range("D" + $row).Value = "Left"
    $tr.Range("E" + $row).Value = "LTR"
    $tr.Range("F" + $row).Value = "PARAMETER"
    $row = $row + 1
    $singleUseId = $singleUseId + 1

    $tr.Range("B" + $row).Value = "SingleUseId" + $singleUseId
    $tr.Range("C" + $row).Value = "Typography_02"
    $tr.Range("D" + $row).Value = "Left"
    $tr.Range("E" + $row).Value = "LTR"
    $tr.Range("F" + $row).Value = "val"
    $row = $row + 1
    $singleUseId = $singleUseId + 1
}
